$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 2120.0527
$ws.Cells.Item(28, 9).Value = 2474.0667
$ws.Cells.Item(28, 10).Value = 792.5
$ws.Cells.Item(28, 11).Value = 2474.0667
$ws.Cells.Item(28, 12).Value = 792.5
$ws.Cells.Item(28, 13).Value = -1989.0667
$ws.Cells.Item(28, 14).Value = -1762.5
$ws.Cells.Item(33, 8).Value = 219.41176
$ws.Cells.Item(33, 9).Value = 123.38461
$ws.Cells.Item(33, 11).Value = 123.38461
$ws.Cells.Item(33, 13).Value = 105.61539
$ws.Cells.Item(51, 8).Value = 1811.8235
$ws.Cells.Item(51, 9).Value = 1100.2
$ws.Cells.Item(51, 10).Value = 2108.3333
$ws.Cells.Item(51, 11).Value = 1100.2
$ws.Cells.Item(51, 12).Value = 2108.3333
$ws.Cells.Item(51, 13).Value = -616.2
$ws.Cells.Item(51, 14).Value = -3076.3333
$ws.Cells.Item(106, 8).Value = 13417.272
$ws.Cells.Item(106, 9).Value = 14158.4
$ws.Cells.Item(106, 11).Value = 14158.4
$ws.Cells.Item(106, 13).Value = -13527.4
$ws.Cells.Item(111, 8).Value = 2803.3572
$ws.Cells.Item(111, 9).Value = 2415.625
$ws.Cells.Item(111, 10).Value = 3320.3333
$ws.Cells.Item(111, 11).Value = 7246.875
$ws.Cells.Item(111, 12).Value = 9960.999899999999
$ws.Cells.Item(111, 13).Value = -4179.875
$ws.Cells.Item(111, 14).Value = -16094.9999
$ws.Cells.Item(113, 8).Value = 2838.5833
$ws.Cells.Item(113, 9).Value = 2696.0
$ws.Cells.Item(113, 10).Value = 3123.75
$ws.Cells.Item(113, 11).Value = 2696.0
$ws.Cells.Item(113, 12).Value = 3123.75
$ws.Cells.Item(113, 13).Value = 558.0
$ws.Cells.Item(113, 14).Value = -9631.75
$ws.Cells.Item(116, 8).Value = 2032.1875
$ws.Cells.Item(116, 10).Value = 2584.3333
$ws.Cells.Item(116, 12).Value = 2584.3333
$ws.Cells.Item(116, 14).Value = -9468.3333
$ws.Cells.Item(118, 8).Value = 640.6923
$ws.Cells.Item(118, 9).Value = 370.0
$ws.Cells.Item(118, 11).Value = 1110.0
$ws.Cells.Item(118, 13).Value = 547.0
$ws.Cells.Item(132, 8).Value = 5559668.5
$ws.Cells.Item(132, 9).Value = 7095245.0
$ws.Cells.Item(132, 10).Value = 7970.077
$ws.Cells.Item(132, 11).Value = 21285735.0
$ws.Cells.Item(132, 12).Value = 23910.231
$ws.Cells.Item(132, 13).Value = -21283205.0
$ws.Cells.Item(132, 14).Value = -28970.231
$ws.Cells.Item(137, 8).Value = 1044.4415
$ws.Cells.Item(137, 9).Value = 814.04
$ws.Cells.Item(137, 11).Value = 2442.12
$ws.Cells.Item(137, 13).Value = 107.8800000000001
$ws.Cells.Item(138, 8).Value = 580894.9
$ws.Cells.Item(138, 9).Value = 664.88574
$ws.Cells.Item(138, 10).Value = 1088596.1
$ws.Cells.Item(138, 11).Value = 1994.65722
$ws.Cells.Item(138, 12).Value = 3265788.3
$ws.Cells.Item(138, 13).Value = 3145.34278
$ws.Cells.Item(138, 14).Value = -3276068.3
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 21323.8
$ws.Cells.Item(2, 9).Value = 1510.5
$ws.Cells.Item(2, 11).Value = 1510.5
$ws.Cells.Item(2, 13).Value = -1397.5
$ws.Cells.Item(32, 8).Value = 4225.058
$ws.Cells.Item(32, 9).Value = 3828.8032
$ws.Cells.Item(32, 11).Value = 3828.8032
$ws.Cells.Item(32, 13).Value = -3541.8032
$ws.Cells.Item(74, 8).Value = 1308.4
$ws.Cells.Item(74, 9).Value = 929.0417
$ws.Cells.Item(74, 10).Value = 2136.0908
$ws.Cells.Item(74, 11).Value = 929.0417
$ws.Cells.Item(74, 12).Value = 2136.0908
$ws.Cells.Item(74, 13).Value = -55.04169999999999
$ws.Cells.Item(74, 14).Value = -3884.0908
$ws.Cells.Item(77, 8).Value = 1308.4
$ws.Cells.Item(77, 9).Value = 929.0417
$ws.Cells.Item(77, 10).Value = 2136.0908
$ws.Cells.Item(77, 11).Value = 4645.2085
$ws.Cells.Item(77, 12).Value = 10680.454
$ws.Cells.Item(77, 13).Value = -277.2084999999997
$ws.Cells.Item(77, 14).Value = -19416.454
$ws.Cells.Item(103, 8).Value = 65000.0
$ws.Cells.Item(103, 10).Value = 65000.0
$ws.Cells.Item(103, 12).Value = 65000.0
$ws.Cells.Item(103, 14).Value = -67344.0
$ws.Cells.Item(116, 8).Value = 21323.8
$ws.Cells.Item(116, 9).Value = 1510.5
$ws.Cells.Item(116, 11).Value = 1510.5
$ws.Cells.Item(116, 13).Value = 783.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 21323.8
$ws.Cells.Item(3, 9).Value = 1510.5
$ws.Cells.Item(3, 11).Value = 1510.5
$ws.Cells.Item(3, 13).Value = -1396.5
$ws.Cells.Item(86, 8).Value = 4491.143
$ws.Cells.Item(86, 9).Value = 5427.6
$ws.Cells.Item(86, 10).Value = 2150.0
$ws.Cells.Item(86, 11).Value = 5427.6
$ws.Cells.Item(86, 12).Value = 2150.0
$ws.Cells.Item(86, 13).Value = -4304.6
$ws.Cells.Item(86, 14).Value = -4396.0
$ws.Cells.Item(89, 8).Value = 4491.143
$ws.Cells.Item(89, 9).Value = 5427.6
$ws.Cells.Item(89, 10).Value = 2150.0
$ws.Cells.Item(89, 11).Value = 27138.0
$ws.Cells.Item(89, 12).Value = 10750.0
$ws.Cells.Item(89, 13).Value = -21522.0
$ws.Cells.Item(89, 14).Value = -21982.0
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 70559.3
$ws.Cells.Item(22, 9).Value = 910.3333
$ws.Cells.Item(22, 11).Value = 910.3333
$ws.Cells.Item(22, 13).Value = -560.3333
$ws.Cells.Item(31, 8).Value = 1934.931
$ws.Cells.Item(31, 9).Value = 2034.7391
$ws.Cells.Item(31, 11).Value = 2034.7391
$ws.Cells.Item(31, 13).Value = -1739.7391
$ws.Cells.Item(34, 8).Value = 1934.931
$ws.Cells.Item(34, 9).Value = 2034.7391
$ws.Cells.Item(34, 11).Value = 2034.7391
$ws.Cells.Item(34, 13).Value = -1832.7391
$ws.Cells.Item(58, 8).Value = 983.5217
$ws.Cells.Item(58, 9).Value = 830.5263
$ws.Cells.Item(58, 11).Value = 830.5263
$ws.Cells.Item(58, 13).Value = -627.5263
$ws.Cells.Item(111, 8).Value = 41266.332
$ws.Cells.Item(111, 10).Value = 41266.332
$ws.Cells.Item(111, 12).Value = 41266.332
$ws.Cells.Item(111, 14).Value = -49446.332
$ws.Cells.Item(114, 8).Value = 25995.0
$ws.Cells.Item(114, 9).Value = 21000.0
$ws.Cells.Item(114, 10).Value = 27660.0
$ws.Cells.Item(114, 11).Value = 21000.0
$ws.Cells.Item(114, 12).Value = 27660.0
$ws.Cells.Item(114, 13).Value = -16661.0
$ws.Cells.Item(114, 14).Value = -36338.0
$ws.Cells.Item(134, 9).Value = 918.9091
$ws.Cells.Item(134, 10).Value = 71430030.0
$ws.Cells.Item(134, 11).Value = 2756.7273
$ws.Cells.Item(134, 12).Value = 214290090.0
$ws.Cells.Item(134, 13).Value = -221.7273
$ws.Cells.Item(134, 14).Value = -214295160.0
$ws.Cells.Item(136, 8).Value = 983.5217
$ws.Cells.Item(136, 9).Value = 830.5263
$ws.Cells.Item(136, 11).Value = 2491.5789
$ws.Cells.Item(136, 13).Value = 58.42110000000002
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 2974.8
$ws.Cells.Item(94, 9).Value = 2249.3333
$ws.Cells.Item(94, 10).Value = 3285.7144
$ws.Cells.Item(94, 11).Value = 6747.999899999999
$ws.Cells.Item(94, 12).Value = 9857.143199999999
$ws.Cells.Item(94, 13).Value = -6071.999899999999
$ws.Cells.Item(94, 14).Value = -11209.1432
$ws.Cells.Item(131, 8).Value = 22728626.0
$ws.Cells.Item(131, 10).Value = 1491.3684
$ws.Cells.Item(131, 12).Value = 4474.1052
$ws.Cells.Item(131, 14).Value = -14554.1052
$ws.Cells.Item(140, 8).Value = 22326.0
$ws.Cells.Item(140, 9).Value = 57766.5
$ws.Cells.Item(140, 10).Value = 2994.818
$ws.Cells.Item(140, 11).Value = 173299.5
$ws.Cells.Item(140, 12).Value = 8984.454000000002
$ws.Cells.Item(140, 13).Value = -168119.5
$ws.Cells.Item(140, 14).Value = -19344.454
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 8).Value = 4000.0
$ws.Cells.Item(47, 10).Value = 4000.0
$ws.Cells.Item(47, 12).Value = 4000.0
$ws.Cells.Item(47, 14).Value = -5136.0
$ws.Cells.Item(132, 8).Value = 2811.4546
$ws.Cells.Item(132, 9).Value = 2215.111
$ws.Cells.Item(132, 10).Value = 5495.0
$ws.Cells.Item(132, 11).Value = 6645.333
$ws.Cells.Item(132, 12).Value = 16485.0
$ws.Cells.Item(132, 13).Value = -4115.333
$ws.Cells.Item(132, 14).Value = -21545.0
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 789.05884
$ws.Cells.Item(22, 9).Value = 528.8889
$ws.Cells.Item(22, 10).Value = 1081.75
$ws.Cells.Item(22, 11).Value = 528.8889
$ws.Cells.Item(22, 12).Value = 1081.75
$ws.Cells.Item(22, 13).Value = -233.8889
$ws.Cells.Item(22, 14).Value = -1671.75
$ws.Cells.Item(27, 8).Value = 789.05884
$ws.Cells.Item(27, 9).Value = 528.8889
$ws.Cells.Item(27, 10).Value = 1081.75
$ws.Cells.Item(27, 11).Value = 528.8889
$ws.Cells.Item(27, 12).Value = 1081.75
$ws.Cells.Item(27, 13).Value = -421.8889
$ws.Cells.Item(27, 14).Value = -1295.75
$ws.Cells.Item(46, 8).Value = 4406.5835
$ws.Cells.Item(46, 9).Value = 1293.0
$ws.Cells.Item(46, 11).Value = 1293.0
$ws.Cells.Item(46, 13).Value = -1105.0
$ws.Cells.Item(55, 8).Value = 328.4643
$ws.Cells.Item(55, 10).Value = 524.5
$ws.Cells.Item(55, 12).Value = 524.5
$ws.Cells.Item(55, 14).Value = -870.5
$ws.Cells.Item(93, 8).Value = 0.0
$ws.Cells.Item(93, 9).Value = 0.0
$ws.Cells.Item(93, 11).Value = 0.0
$ws.Cells.Item(93, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 19731.527
$ws.Cells.Item(132, 9).Value = 1203.5
$ws.Cells.Item(132, 10).Value = 49729.285
$ws.Cells.Item(132, 11).Value = 3610.5
$ws.Cells.Item(132, 12).Value = 149187.855
$ws.Cells.Item(132, 13).Value = -1080.5
$ws.Cells.Item(132, 14).Value = -154247.855
$ws.Cells.Item(136, 8).Value = 1105.359
$ws.Cells.Item(136, 9).Value = 988.5
$ws.Cells.Item(136, 10).Value = 1900.0
$ws.Cells.Item(136, 11).Value = 2965.5
$ws.Cells.Item(136, 12).Value = 5700.0
$ws.Cells.Item(136, 13).Value = -415.5
$ws.Cells.Item(136, 14).Value = -10800.0
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 4057.1428
$ws.Cells.Item(96, 9).Value = 2550.0
$ws.Cells.Item(96, 11).Value = 2550.0
$ws.Cells.Item(96, 13).Value = -1177.0
$ws.Cells.Item(107, 8).Value = 471.1875
$ws.Cells.Item(107, 9).Value = 431.66666
$ws.Cells.Item(107, 11).Value = 1294.99998
$ws.Cells.Item(107, 13).Value = 625.00002
$ws.Cells.Item(136, 8).Value = 781.0417
$ws.Cells.Item(136, 9).Value = 524.6
$ws.Cells.Item(136, 10).Value = 1208.4445
$ws.Cells.Item(136, 11).Value = 1573.8
$ws.Cells.Item(136, 12).Value = 3625.3335
$ws.Cells.Item(136, 13).Value = 976.1999999999998
$ws.Cells.Item(136, 14).Value = -8725.3335
